# Restore C10 ("R30" row, "C1" / min<=hour column) from 18 to 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C10").Value = 1
